$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the metadata "type" row (row 2): dimensions curated -> some become measures,
# and aragon becomes an sdmx refArea dimension.
$ws.Range("B2").Value = "iaest-measure:sector-actividad"
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "iaest-measure:sexo"

# Update the "medida/dim" classification row (row 3) to match the new measure columns.
$ws.Range("B3").Value = "medida"
$ws.Range("F3").Value = "medida"

# Update the data-type row (row 4): former skos:Concept columns are now plain xsd:int,
# except the refArea column which now references a URI-Comunidad type.
$ws.Range("B4").Value = "xsd:int"
$ws.Range("E4").Value = "URI-Comunidad"
$ws.Range("F4").Value = "xsd:int"

# Remove row 5 entirely (the old mapping-*.xlsx reference row is no longer needed).
$ws.Rows.Item(5).Delete()
